$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "298.77"
Set-TextValue "E2" "-2.02%"
Set-TextValue "G2" "3"
Set-TextValue "D3" "31.30"
Set-TextValue "G3" "3"
Set-TextValue "D4" "5.090"
Set-TextValue "E4" "-2.00%"
Set-TextValue "G4" "3"
Set-TextValue "D5" "0.07899"
Set-TextValue "E5" "5.91%"
Set-TextValue "G5" "3"
Set-TextValue "D6" "2.339"
Set-TextValue "E6" "-3.40%"
Set-TextValue "G6" "3"
Set-TextValue "D7" "7.739"
Set-TextValue "E7" "-3.09%"
Set-TextValue "G7" "3"
Set-TextValue "D8" "3.864"
Set-TextValue "E8" "0.01%"
Set-TextValue "G8" "3"
Set-TextValue "D9" "0.9200"
Set-TextValue "E9" "-0.05%"
Set-TextValue "G9" "3"
Set-TextValue "D10" "0.1732"
Set-TextValue "E10" "0.09%"
Set-TextValue "G10" "3"
Set-TextValue "D11" "0.07326"
Set-TextValue "E11" "-4.93%"
Set-TextValue "G11" "3"
Set-TextValue "D12" "0.09005"
Set-TextValue "E12" "9.65%"
Set-TextValue "G12" "3"
Set-TextValue "D13" "0.03050"
Set-TextValue "E13" "1.34%"
Set-TextValue "G13" "3"
Set-TextValue "D14" "0.1002"
Set-TextValue "E14" "0.88%"
Set-TextValue "G14" "3"
Set-TextValue "D15" "0.001512"
Set-TextValue "E15" "0.09%"
Set-TextValue "G15" "3"
Set-TextValue "D16" "0.006065"
Set-TextValue "E16" "-0.83%"
Set-TextValue "G16" "3"
Set-TextValue "D17" "3.479"
Set-TextValue "E17" "-0.41%"
Set-TextValue "G17" "3"
Set-TextValue "D18" "2.266"
Set-TextValue "E18" "1.67%"
Set-TextValue "G18" "3"
Set-TextValue "E19" "0.30%"
Set-TextValue "G19" "3"
Set-TextValue "D20" "0.1318"
Set-TextValue "E20" "-1.34%"
Set-TextValue "G20" "3"
Set-TextValue "D21" "4.161"
Set-TextValue "E21" "-10.60%"
Set-TextValue "G21" "3"
Set-TextValue "D22" "0.1699"
Set-TextValue "E22" "8.56%"
Set-TextValue "G22" "3"
Set-TextValue "D23" "0.04616"
Set-TextValue "E23" "0.07%"
Set-TextValue "G23" "3"
Set-TextValue "D24" "0.001241"
Set-TextValue "E24" "-1.17%"
Set-TextValue "G24" "3"
Set-TextValue "E25" "-1.31%"
Set-TextValue "G25" "3"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "-7.65%"
Set-TextValue "G26" "3"
Set-TextValue "D27" "0.0003395"
Set-TextValue "E27" "23.92%"
Set-TextValue "G27" "3"
Set-TextValue "G28" "3"
Set-TextValue "G29" "3"
Set-TextValue "G30" "3"
Set-TextValue "G31" "3"
Set-TextValue "G32" "3"
Set-TextValue "G33" "3"
Set-TextValue "G34" "3"
Set-TextValue "G35" "3"
Set-TextValue "G36" "3"
Set-TextValue "G37" "3"
Set-TextValue "G38" "3"
Set-TextValue "D39" "0.01741"
Set-TextValue "E39" "-1.89%"
Set-TextValue "G39" "3"
Set-TextValue "D40" "0.04600"
Set-TextValue "E40" "0.97%"
Set-TextValue "G40" "3"
Set-TextValue "D41" "0.006982"
Set-TextValue "E41" "-5.29%"
Set-TextValue "G41" "3"
Set-TextValue "G42" "3"
Set-TextValue "D43" "0.002189"
Set-TextValue "E43" "0.50%"
Set-TextValue "G43" "3"
Set-TextValue "D44" "0.009539"
Set-TextValue "E44" "-11.70%"
Set-TextValue "G44" "3"
Set-TextValue "D45" "0.00006281"
Set-TextValue "E45" "2.80%"
Set-TextValue "G45" "3"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "E46" "-0.03%"
Set-TextValue "G46" "3"
Set-TextValue "D47" "0.007974"
Set-TextValue "E47" "-19.38%"
Set-TextValue "G47" "3"
Set-TextValue "D48" "0.7475"
Set-TextValue "E48" "-8.90%"
Set-TextValue "G48" "3"
Set-TextValue "D49" "0.00002099"
Set-TextValue "E49" "-0.03%"
Set-TextValue "G49" "3"
Set-TextValue "D50" "0.0001999"
Set-TextValue "E50" "0.04%"
Set-TextValue "G50" "3"
Set-TextValue "G51" "3"
